$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BENCHMARK")

# Row 2: E2 cleared
$ws.Range("E2").Value = ""

# Column E filled in for rows 3-6, 8-11 (matching neighboring values in C/D/F)
$ws.Range("E3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("E4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("E5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("E6").Value = "6,09 TL - 12,19 TL - 152,35 TL"

$ws.Range("E8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("E9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("E10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("E11").Value = "3,04 TL - 6,09 TL - 76,17 TL"

# Row 13: C13 updated, E13 filled in
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 0,94 TL"
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 851,5 TL"

# Row 14: E14 filled in
$ws.Range("E14").Value = "1.660 TL - 1.660 TL"
